# Update the date line and the 25 division-fact answers in the table.
# Text-based Find/Replace (whole-document, whole-match, replace-all) is
# safe here because every "old" string in the document is unique EXCEPT
# for one collision: the new value destined for one cell
# ("186÷4=46, 2") equals another cell's *old* value. We replace that
# original occurrence first (near the top of the script) so that later,
# when "186÷4=46, 2" is written in as a *new* value, it is the only
# instance of that text left in the document.

$d = $word.ActiveDocument

function Replace-Text($old, $new, $count) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $new, $count)
}

# --- Date heading -----------------------------------------------------
Replace-Text "2025-02-17 Monday" "2025-02-18 Tuesday" 2

# --- Collision pair: the NEW text of cell #24 equals the OLD text of
#     cell #6 ("186÷4=46, 2"). Retire the original occurrence FIRST, so
#     that when cell #24 is later rewritten to the same string there is
#     only one (correct) instance of it left in the document.
Replace-Text "186÷4=46, 2" "427÷4=106, 3" 2

# --- Remaining table cells (no collisions) -----------------------------
Replace-Text "591÷9=65, 6" "990÷8=123, 6" 2
Replace-Text "696÷4=174, 0" "214÷3=71, 1" 2
Replace-Text "433÷2=216, 1" "368÷9=40, 8" 2
Replace-Text "486÷9=54, 0" "512÷2=256, 0" 2
Replace-Text "573÷7=81, 6" "708÷2=354, 0" 2

Replace-Text "622÷7=88, 6" "397÷4=99, 1" 2
Replace-Text "620÷3=206, 2" "840÷3=280, 0" 2
Replace-Text "596÷3=198, 2" "511÷7=73, 0" 2
Replace-Text "756÷9=84, 0" "240÷6=40, 0" 2

Replace-Text "637÷7=91, 0" "372÷7=53, 1" 2
Replace-Text "837÷2=418, 1" "548÷8=68, 4" 2
Replace-Text "668÷4=167, 0" "538÷7=76, 6" 2
Replace-Text "342÷8=42, 6" "810÷8=101, 2" 2
Replace-Text "666÷2=333, 0" "224÷5=44, 4" 2

Replace-Text "523÷2=261, 1" "366÷9=40, 6" 2
Replace-Text "296÷8=37, 0" "236÷9=26, 2" 2
Replace-Text "576÷2=288, 0" "847÷2=423, 1" 2
Replace-Text "824÷3=274, 2" "194÷2=97, 0" 2
Replace-Text "512÷9=56, 8" "439÷4=109, 3" 2

Replace-Text "474÷6=79, 0" "776÷7=110, 6" 2
Replace-Text "420÷8=52, 4" "562÷4=140, 2" 2
Replace-Text "749÷2=374, 1" "318÷7=45, 3" 2
Replace-Text "382÷4=95, 2" "186÷4=46, 2" 2
Replace-Text "909÷5=181, 4" "763÷5=152, 3" 2
